$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 150
$ws.Range("I8").Value = 150
$ws.Range("K8").Value = 450
$ws.Range("M8").Value = -311
$ws.Range("H46").Value = 1614.7778
$ws.Range("J46").Value = 839.6667
$ws.Range("L46").Value = 2519.0001
$ws.Range("N46").Value = -2757.0001
$ws.Range("H60").Value = 1614.7778
$ws.Range("J60").Value = 839.6667
$ws.Range("L60").Value = 2519.0001
$ws.Range("N60").Value = -3487.0001
$ws.Range("H113").Value = 3161.8
$ws.Range("I113").Value = 3002.25
$ws.Range("J113").Value = 3800
$ws.Range("K113").Value = 3002.25
$ws.Range("L113").Value = 3800
$ws.Range("M113").Value = 251.75
$ws.Range("N113").Value = -10308
$ws.Range("H129").Value = 1304.3793
$ws.Range("I129").Value = 3398.6
$ws.Range("J129").Value = 868.0833
$ws.Range("K129").Value = 10195.8
$ws.Range("L129").Value = 2604.2499
$ws.Range("M129").Value = -5195.799999999999
$ws.Range("N129").Value = -12604.2499
$ws.Range("H134").Value = 52385.31
$ws.Range("J134").Value = 52385.31
$ws.Range("L134").Value = 52385.31
$ws.Range("N134").Value = -62525.31
$ws.Range("H141").Value = 3704.7
$ws.Range("I141").Value = 2506.3333
$ws.Range("K141").Value = 7518.999899999999
$ws.Range("M141").Value = -2338.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 537.6
$ws.Range("I5").Value = 296
$ws.Range("J5").Value = 900
$ws.Range("K5").Value = 296
$ws.Range("L5").Value = 900
$ws.Range("M5").Value = -184
$ws.Range("N5").Value = -1124
$ws.Range("H32").Value = 16480.303
$ws.Range("I32").Value = 16623.2
$ws.Range("J32").Value = 15855.125
$ws.Range("K32").Value = 16623.2
$ws.Range("L32").Value = 15855.125
$ws.Range("M32").Value = -16336.2
$ws.Range("N32").Value = -16429.125
$ws.Range("H37").Value = 38478.184
$ws.Range("I37").Value = 1000
$ws.Range("J37").Value = 42226
$ws.Range("K37").Value = 1000
$ws.Range("L37").Value = 42226
$ws.Range("M37").Value = -727
$ws.Range("N37").Value = -42772
$ws.Range("H45").Value = 794.3200000000001
$ws.Range("I45").Value = 794.46466
$ws.Range("K45").Value = 794.46466
$ws.Range("M45").Value = -417.46466
$ws.Range("H55").Value = 42999.332
$ws.Range("J55").Value = 42999.332
$ws.Range("L55").Value = 42999.332
$ws.Range("N55").Value = -43629.332
$ws.Range("H110").Value = 4349.9473
$ws.Range("I110").Value = 1896.0769
$ws.Range("K110").Value = 1896.0769
$ws.Range("M110").Value = 148.9231
$ws.Range("H132").Value = 15153539
$ws.Range("I132").Value = 17858602
$ws.Range("J132").Value = 5188.8
$ws.Range("K132").Value = 53575806
$ws.Range("L132").Value = 15566.4
$ws.Range("M132").Value = -53573276
$ws.Range("N132").Value = -20626.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 537.6
$ws.Range("I4").Value = 296
$ws.Range("J4").Value = 900
$ws.Range("K4").Value = 296
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = -181
$ws.Range("N4").Value = -1130
$ws.Range("H107").Value = 3571.6667
$ws.Range("I107").Value = 2766.5
$ws.Range("K107").Value = 2766.5
$ws.Range("M107").Value = -846.5
$ws.Range("H134").Value = 1736.6333
$ws.Range("I134").Value = 1200.125
$ws.Range("K134").Value = 3600.375
$ws.Range("M134").Value = -1065.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 554.1429000000001
$ws.Range("I7").Value = 378
$ws.Range("J7").Value = 1200
$ws.Range("K7").Value = 378
$ws.Range("L7").Value = 1200
$ws.Range("M7").Value = -265
$ws.Range("N7").Value = -1426
$ws.Range("H70").Value = 34526.668
$ws.Range("J70").Value = 34526.668
$ws.Range("L70").Value = 34526.668
$ws.Range("N70").Value = -35156.668
$ws.Range("H73").Value = 34526.668
$ws.Range("J73").Value = 34526.668
$ws.Range("L73").Value = 34526.668
$ws.Range("N73").Value = -36710.668
$ws.Range("H88").Value = 24886.8
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 28608.5
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 28608.5
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -29420.5
$ws.Range("H91").Value = 24886.8
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 28608.5
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 28608.5
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -31416.5
$ws.Range("H112").Value = 40283
$ws.Range("J112").Value = 40283
$ws.Range("L112").Value = 40283
$ws.Range("N112").Value = -43237
$ws.Range("H141").Value = 5403.3335
$ws.Range("J141").Value = 5403.3335
$ws.Range("L141").Value = 5403.3335
$ws.Range("N141").Value = -15763.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1885
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 2480
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 7440
$ws.Range("M4").Value = -188
$ws.Range("N4").Value = -7664
$ws.Range("H11").Value = 406.8889
$ws.Range("I11").Value = 563.75
$ws.Range("J11").Value = 281.4
$ws.Range("K11").Value = 1691.25
$ws.Range("L11").Value = 844.1999999999999
$ws.Range("M11").Value = -1551.25
$ws.Range("N11").Value = -1124.2
$ws.Range("H33").Value = 20534692
$ws.Range("I33").Value = 94.5
$ws.Range("J33").Value = 24268254
$ws.Range("K33").Value = 567
$ws.Range("L33").Value = 145609524
$ws.Range("M33").Value = -284
$ws.Range("N33").Value = -145610090
$ws.Range("H131").Value = 937.4400000000001
$ws.Range("I131").Value = 812.1667
$ws.Range("J131").Value = 945.43616
$ws.Range("K131").Value = 2436.5001
$ws.Range("L131").Value = 2836.30848
$ws.Range("M131").Value = 2603.4999
$ws.Range("N131").Value = -12916.30848
$ws.Range("H139").Value = 152680.16
$ws.Range("I139").Value = 202126.67
$ws.Range("J139").Value = 4340.6
$ws.Range("K139").Value = 606380.01
$ws.Range("L139").Value = 13021.8
$ws.Range("M139").Value = -601240.01
$ws.Range("N139").Value = -23301.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2308.3845
$ws.Range("I102").Value = 1641.8
$ws.Range("J102").Value = 2725
$ws.Range("K102").Value = 1641.8
$ws.Range("L102").Value = 2725
$ws.Range("M102").Value = -19.79999999999995
$ws.Range("N102").Value = -5969
$ws.Range("H122").Value = 2799.0625
$ws.Range("I122").Value = 2405.3076
$ws.Range("J122").Value = 4505.3335
$ws.Range("K122").Value = 7215.9228
$ws.Range("L122").Value = 13516.0005
$ws.Range("M122").Value = -4765.9228
$ws.Range("N122").Value = -18416.0005
$ws.Range("H126").Value = 16125.19
$ws.Range("I126").Value = 34904.332
$ws.Range("J126").Value = 2040.8334
$ws.Range("K126").Value = 104712.996
$ws.Range("L126").Value = 6122.5002
$ws.Range("M126").Value = -102242.996
$ws.Range("N126").Value = -11062.5002
$ws.Range("H132").Value = 2933.3333
$ws.Range("I132").Value = 2272.0715
$ws.Range("J132").Value = 5247.75
$ws.Range("K132").Value = 6816.2145
$ws.Range("L132").Value = 15743.25
$ws.Range("M132").Value = -4286.2145
$ws.Range("N132").Value = -20803.25
$ws.Range("H140").Value = 37263.6
$ws.Range("J140").Value = 37263.6
$ws.Range("L140").Value = 37263.6
$ws.Range("N140").Value = -47623.6
$ws.Range("H141").Value = 70103.75
$ws.Range("J141").Value = 70103.75
$ws.Range("L141").Value = 70103.75
$ws.Range("N141").Value = -80463.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1342.7142
$ws.Range("I22").Value = 1199.6666
$ws.Range("J22").Value = 1450
$ws.Range("K22").Value = 1199.6666
$ws.Range("L22").Value = 1450
$ws.Range("M22").Value = -904.6666
$ws.Range("N22").Value = -2040
$ws.Range("H27").Value = 1342.7142
$ws.Range("I27").Value = 1199.6666
$ws.Range("J27").Value = 1450
$ws.Range("K27").Value = 1199.6666
$ws.Range("L27").Value = 1450
$ws.Range("M27").Value = -1092.6666
$ws.Range("N27").Value = -1664
$ws.Range("H122").Value = 2238.0952
$ws.Range("I122").Value = 2225
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6675
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4225
$ws.Range("N122").Value = -12400
$ws.Range("H138").Value = 53477.5
$ws.Range("J138").Value = 53477.5
$ws.Range("L138").Value = 53477.5
$ws.Range("N138").Value = -63757.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 2000
$ws.Range("I21").Value = 2000
$ws.Range("K21").Value = 2000
$ws.Range("M21").Value = -1765
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1710
$ws.Range("H132").Value = 18520458
$ws.Range("I132").Value = 1585.4375
$ws.Range("K132").Value = 4756.3125
$ws.Range("M132").Value = -2226.3125
$ws.Range("H135").Value = 32452.941
$ws.Range("J135").Value = 32452.941
$ws.Range("L135").Value = 32452.941
$ws.Range("N135").Value = -42592.941
$ws.Range("H140").Value = 34966.77
$ws.Range("J140").Value = 34966.77
$ws.Range("L140").Value = 34966.77
$ws.Range("N140").Value = -45326.77
$ws.Range("H141").Value = 38578.855
$ws.Range("J141").Value = 38578.855
$ws.Range("L141").Value = 38578.855
$ws.Range("N141").Value = -49326.77
